# tcm-words.xlsx maintenance edit:
#  - "症状" (symptoms) sheet had a duplicated row (row 10 was an exact
#    duplicate of row 9: "肥胖" / "obesity") left over from a copy/paste
#    before the word-trim fix landed. Remove the stray duplicate row so
#    the list is clean again (rows below shift up, used range shrinks).
#  - Update the saved cursor/selection state on both affected sheets to
#    match where the editor was left after making the change.

$wb = $excel.ActiveWorkbook

$wsZhongyi = $wb.Worksheets.Item(1)   # 中医 (tab stays active throughout)
$wsZhengzhuang = $wb.Worksheets.Item(3)   # 症状

# Remove the duplicate "肥胖 / obesity" row (row 10) from 症状.
# EntireRow delete shifts rows 11:17 up to 10:16 and shrinks the sheet's
# used range from A1:B17 down to A1:B16 automatically.
$wsZhengzhuang.Rows(10).Delete()

# Leave the editor's selection on 症状 where it was parked afterwards.
$wsZhengzhuang.Activate()
$wsZhengzhuang.Range("B22").Select()

# Re-activate 中医 (it was the selected tab before/after the edit) and
# restore its parked selection too.
$wsZhongyi.Activate()
$wsZhongyi.Range("D1:E4").Select()
